$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "Daily.jamais.frqce_conso"
$ws.Range("C1").Value = "4 to 6 days per week.jamais.frqce_conso"
$ws.Range("D1").Value = "2 to 3 days per week.jamais.frqce_conso"
$ws.Range("E1").Value = "Once a week or less.jamais.frqce_conso"
$ws.Range("F1").Value = "Not used in the last 30 days.jamais.frqce_conso"
$ws.Range("G1").Value = "Not known / missing.jamais.frqce_conso"
$ws.Range("H1").Value = "Total.jamais.frqce_conso"
